$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Shreya Sahni's Feb 1st week progress update
$ws.Range("B6").Value = "1. Studied syntax and other basics of dart, flutter 2. Set up flutter plugin with android studio"

$ws.Range("B14").Select()
